{"js": "// Applies the resume \"abstraction\" edits described in the diff:\n//  1. Name -> title case\n//  2. Title -> generic placeholder\n//  3. Phone/email formatting -> generic/lowercase formatting\n//  4. \"20+ years\" -> \"21 years\"\n//  5. Company name -> generic placeholder\n//  6. Remove product names (BALLISTA and DAMON) from a bullet\n//  7. Remove five entire job entries (DATA PRODUCTS MANAGER ... RESEARCH\n//     DIRECTOR & PRODUCT MANAGER) from PROFESSIONAL EXPERIENCE\n//  8. Remove product names from the achievements bullet\n//  9. Remove the two trailing achievement sub-sections (Market\n//     Intelligence & Research Leadership; Cross-Functional Leadership &\n//     Collaboration)\n\nconst body = context.document.body;\n\n// --- 1. Name -------------------------------------------------------------\nconst nameResults = body.search(\"DHEERAJ CHAND\", { matchCase: true });\nnameResults.load(\"items\");\nawait context.sync();\nif (nameResults.items.length > 0) {\n  nameResults.items[0].insertText(\"Dheeraj Chand\", Word.InsertLocation.replace);\n}\n\n// --- 2. Title --------------------------------------------------------------\nconst titleResults = body.search(\"Senior Product Marketing Manager\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(\"Professional Title\", Word.InsertLocation.replace);\n}\n\n// --- 3. Contact line ---------------------------------------------------\nconst contactResults = body.search(\"(202) 550-7110 | Dheeraj.Chand@gmail.com\", { matchCase: true });\ncontactResults.load(\"items\");\nawait context.sync();\nif (contactResults.items.length > 0) {\n  contactResults.items[0].insertText(\"202.550.7110 | dheeraj.chand@gmail.com\", Word.InsertLocation.replace);\n}\n\n// --- 4. Years of experience ----------------------------------------------\nconst yearsResults = body.search(\"20+ years\", { matchCase: true });\nyearsResults.load(\"items\");\nawait context.sync();\nif (yearsResults.items.length > 0) {\n  yearsResults.items[0].insertText(\"21 years\", Word.InsertLocation.replace);\n}\n\n// --- 5. Company name (Siege Analytics line) -------------------------------\nconst companyResults = body.search(\"Siege Analytics, Austin, TX | 2005 \u2013 Present\", { matchCase: true });\ncompanyResults.load(\"items\");\nawait context.sync();\nif (companyResults.items.length > 0) {\n  companyResults.items[0].insertText(\"Your Company Name, Your City, ST | 2005 \u2013 Present\", Word.InsertLocation.replace);\n}\n\n// --- 6. Remove product names from the go-to-market bullet ----------------\nconst gtmResults = body.search(\n  \"Developed and executed go-to-market strategies for multiple SaaS platform launches including BALLISTA and DAMON, achieving thousands of active users and significant market penetration\",\n  { matchCase: true }\n);\ngtmResults.load(\"items\");\nawait context.sync();\nif (gtmResults.items.length > 0) {\n  gtmResults.items[0].insertText(\n    \"Developed and executed go-to-market strategies for multiple SaaS platform launches, achieving thousands of active users and significant market penetration\",\n    Word.InsertLocation.replace\n  );\n}\n\n// --- 7. Remove the five whole job entries ---------------------------------\n// They run from the \"DATA PRODUCTS MANAGER\" heading through the last bullet\n// of \"RESEARCH DIRECTOR & PRODUCT MANAGER\" (just before \"KEY ACHIEVEMENTS\n// AND IMPACT\").\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nlet startIdx = -1;\nlet endIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === \"DATA PRODUCTS MANAGER\") {\n    startIdx = i;\n  }\n  if (items[i].text === \"KEY ACHIEVEMENTS AND IMPACT\") {\n    endIdx = i;\n    break;\n  }\n}\n\nif (startIdx !== -1 && endIdx !== -1) {\n  // Delete paragraph objects individually (each holds its own identity, so\n  // deleting from the bottom up is not even required, but we do it anyway\n  // to keep things tidy) \u2014 this removes the paragraph mark cleanly, unlike\n  // deleting a Range spanning multiple paragraphs which can leave an empty\n  // paragraph behind.\n  for (let i = endIdx - 1; i >= startIdx; i--) {\n    items[i].delete();\n  }\n  await context.sync();\n}\n\n// --- 8. Remove product names from the achievements bullet -----------------\nconst achResults = body.search(\n  \"Successfully launched multiple B2B SaaS platforms (BALLISTA, DAMON, SimCrisis, RACSO) used by thousands of active users with proven market adoption and customer retention\",\n  { matchCase: true }\n);\nachResults.load(\"items\");\nawait context.sync();\nif (achResults.items.length > 0) {\n  achResults.items[0].insertText(\n    \"Successfully launched multiple B2B SaaS platforms used by thousands of active users with proven market adoption and customer retention\",\n    Word.InsertLocation.replace\n  );\n}\n\n// --- 9. Remove the two trailing achievement sub-sections -------------------\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"text\");\nawait context.sync();\n\nconst items2 = paragraphs2.items;\nlet tailStart = -1;\nfor (let i = 0; i < items2.length; i++) {\n  if (items2[i].text === \"Market Intelligence & Research Leadership\") {\n    tailStart = i;\n    break;\n  }\n}\n\nif (tailStart !== -1) {\n  const lastIdx = items2.length - 1;\n  for (let i = lastIdx; i >= tailStart; i--) {\n    items2[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Applies the resume \"abstraction\" edits described in the diff:\n#  1. Name -> title case\n#  2. Title -> generic placeholder\n#  3. Phone/email formatting -> generic/lowercase formatting\n#  4. \"20+ years\" -> \"21 years\"\n#  5. Company name -> generic placeholder\n#  6. Remove product names (BALLISTA and DAMON) from a bullet\n#  7. Remove five entire job entries (DATA PRODUCTS MANAGER ... RESEARCH\n#     DIRECTOR & PRODUCT MANAGER) from PROFESSIONAL EXPERIENCE\n#  8. Remove product names from the achievements bullet\n#  9. Remove the two trailing achievement sub-sections (Market\n#     Intelligence & Research Leadership; Cross-Functional Leadership &\n#     Collaboration)\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# --- 1. Name ---------------------------------------------------------------\nReplace-Text \"DHEERAJ CHAND\" \"Dheeraj Chand\"\n\n# --- 2. Title ----------------------------------------------------------------\nReplace-Text \"Senior Product Marketing Manager\" \"Professional Title\"\n\n# --- 3. Contact line ---------------------------------------------------------\nReplace-Text \"(202) 550-7110 | Dheeraj.Chand@gmail.com\" \"202.550.7110 | dheeraj.chand@gmail.com\"\n\n# --- 4. Years of experience ---------------------------------------------------\nReplace-Text \"20+ years\" \"21 years\"\n\n# --- 5. Company name (Siege Analytics line) -----------------------------------\nReplace-Text \"Siege Analytics, Austin, TX | 2005 \u2013 Present\" \"Your Company Name, Your City, ST | 2005 \u2013 Present\"\n\n# --- 6. Remove product names from the go-to-market bullet ---------------------\nReplace-Text `\n  \"Developed and executed go-to-market strategies for multiple SaaS platform launches including BALLISTA and DAMON, achieving thousands of active users and significant market penetration\" `\n  \"Developed and executed go-to-market strategies for multiple SaaS platform launches, achieving thousands of active users and significant market penetration\"\n\n# --- 7. Remove the five whole job entries -------------------------------------\n# They run from the \"DATA PRODUCTS MANAGER\" heading through the last bullet of\n# \"RESEARCH DIRECTOR & PRODUCT MANAGER\" (just before \"KEY ACHIEVEMENTS AND\n# IMPACT\").\n$startIdx = -1\n$endIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\")\n    if ($t -eq \"DATA PRODUCTS MANAGER\") { $startIdx = $i }\n    if ($t -eq \"KEY ACHIEVEMENTS AND IMPACT\") { $endIdx = $i; break }\n}\n\nif ($startIdx -ne -1 -and $endIdx -ne -1) {\n    $startPara = $d.Paragraphs.Item($startIdx)\n    $endPara = $d.Paragraphs.Item($endIdx - 1)\n    $rangeToDelete = $d.Range($startPara.Range.Start, $endPara.Range.End)\n    $rangeToDelete.Delete()\n}\n\n# --- 8. Remove product names from the achievements bullet ----------------------\nReplace-Text `\n  \"Successfully launched multiple B2B SaaS platforms (BALLISTA, DAMON, SimCrisis, RACSO) used by thousands of active users with proven market adoption and customer retention\" `\n  \"Successfully launched multiple B2B SaaS platforms used by thousands of active users with proven market adoption and customer retention\"\n\n# --- 9. Remove the two trailing achievement sub-sections ------------------------\n$tailStart = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\")\n    if ($t -eq \"Market Intelligence & Research Leadership\") { $tailStart = $i; break }\n}\n\nif ($tailStart -ne -1) {\n    $lastIdx = $d.Paragraphs.Count\n    $tailStartPara = $d.Paragraphs.Item($tailStart)\n    $tailEndPara = $d.Paragraphs.Item($lastIdx)\n    $tailRange = $d.Range($tailStartPara.Range.Start, $tailEndPara.Range.End)\n    $tailRange.Delete()\n}\n"}
